# Worked on in and out
# Add three new timesheet entries (rows 21-23) to Sheet1, matching the
# existing date/hour entries above them, then let Excel recalculate the
# dependent formulas (L4, P4, M9, P9) and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the date format already applied to the preceding rows (A17:A20).
$dateFormat = $ws.Range("A20").NumberFormat

$ws.Range("A21").Value = 43935
$ws.Range("A21").NumberFormat = $dateFormat
$ws.Range("B21").Value = 1

$ws.Range("A22").Value = 43936
$ws.Range("A22").NumberFormat = $dateFormat
$ws.Range("B22").Value = 2

$ws.Range("A23").Value = 43938
$ws.Range("A23").NumberFormat = $dateFormat
$ws.Range("B23").Value = 3

# Recalculate dependent formulas (SUM, totals, owed) now that new rows exist.
$excel.Calculate()

# Match the author's final cursor/selection position.
$ws.Range("D22").Select()
